$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings get appended in the exact order their cell .Value is
# set, so the order below is chosen to reproduce the target sharedStrings
# table order (E12, D13, E13, C14, D14, E14, C12, D12).

# Row 12 (test case 6): actual results
$ws.Cells.Item(12, 5).Value = "los cambios se realizan y persisten sin problemas"

# Row 13 (test case 7): expected / actual results
$ws.Cells.Item(13, 4).Value = "Poder eliminar el registro seleccionado"
$ws.Cells.Item(13, 5).Value = "El registro se elimina sin problemas"

# Row 14 (test case 8): action / expected / actual results
$ws.Cells.Item(14, 3).Value = "Buscar con cualquier dato de un Host"
$ws.Cells.Item(14, 4).Value = "Encontrar un registro mediante el ingreso de un valor de referencia"
$ws.Cells.Item(14, 5).Value = "Devuelve correctamente resultados si hay coincidencia"

# Row 12 (test case 6): action / expected results
$ws.Cells.Item(12, 3).Value = "Editar los campos correspondientes en contactos"
$ws.Cells.Item(12, 4).Value = "cambiar cualquiera de los campos correspondiente a un contacto"

# Row heights
$ws.Rows.Item(12).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 45

# Update selection to match new active cell
$ws.Range("D16").Select()
